# Add a new "2022" column (S) to the right of the existing "2021" column (R),
# mirroring the formatting of column R for rows 3-34 and filling in the new
# year's data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the formatting (styles/borders/number formats) from R3:R34 into
#    the new S3:S34 range so the new column visually matches its neighbour.
$ws.Range("R3:R34").Copy() | Out-Null
$ws.Range("S3:S34").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2. Header (year) value for the new column.
$ws.Range("S4").Value = 2022

# 3. Data values for the new column, row by row.
$ws.Range("S5").Value = 135
$ws.Range("S6").Value = 99
$ws.Range("S7").Value = 36
$ws.Range("S8").Value = 97
$ws.Range("S9").Value = 80
$ws.Range("S10").Value = 17
$ws.Range("S11").Value = 17
$ws.Range("S12").Value = 11
$ws.Range("S13").Value = 6
$ws.Range("S14").Value = 5
$ws.Range("S15").Value = 3
$ws.Range("S16").Value = 2
$ws.Range("S17").Value = "-"
$ws.Range("S18").Value = "-"
$ws.Range("S19").Value = "-"
$ws.Range("S20").Value = 6
$ws.Range("S21").Value = 1
$ws.Range("S22").Value = 5
$ws.Range("S23").Value = "-"
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = "-"
$ws.Range("S26").Value = 10
$ws.Range("S27").Value = 4
$ws.Range("S28").Value = 6
$ws.Range("S29").Value = "-"
$ws.Range("S30").Value = "-"
$ws.Range("S31").Value = "-"
$ws.Range("S32").Value = "-"
$ws.Range("S33").Value = "-"
$ws.Range("S34").Value = "-"

# 4. Move the selection where the author left it after entering the data
#    (one column to the right, at the header row).
$ws.Range("T4").Select()
